$d = $word.ActiveDocument

# --- Add the three new character styles (matching the OOXML <w:style> defs) ---

$sStyle = $d.Styles.Add("GaNStyle", 2)
$sStyle.Font.Name = "Calibri"
$sStyle.Font.NameAscii = "Calibri"
$sStyle.Font.Size = 14

$sParagraph = $d.Styles.Add("GaNParagraph", 2)
$sParagraph.Font.Name = "Calibri"
$sParagraph.Font.NameAscii = "Calibri"
$sParagraph.Font.Size = 10

$sLinks = $d.Styles.Add("GaNLinks", 2)
$sLinks.Font.Name = "Calibri"
$sLinks.Font.NameAscii = "Calibri"
$sLinks.Font.Bold = $true
$sLinks.Font.Color = 8388608
$sLinks.Font.Size = 9.5
$sLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Datas da campaña..." run (4 occurrences) ---

$datesText = "Datas da campaña de 2022 que usan Constelación de Géminis: 14-23 de febreiro, 14-24 de marzo"
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($datesText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $range.Style = "GaNStyle"
    $range.Collapse(0)
    $found = $range.Find.Execute($datesText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNParagraph to the "Grazas por participar..." run ---

$thanksText = "Grazas por participar nesta campaña global de medida da contaminación lumínica mediante a observación das estrelas máis febles que podes albiscar. Localizando e observando a  Constelación de Géminis e comparándoa co que aparece nos mapas estelares recollidos neste documento podes saber canto contribúen á contaminación lumínica os sistemas de iluminación que hai no teu barrio ou vila. As túas achegas á base de datos en liña de GLOBE at Night (O MUNDO á Noite) servirán para documentar a calidade do ceo nocturno."
$range2 = $d.Content
$range2.Find.ClearFormatting()
$found2 = $range2.Find.Execute($thanksText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $range2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Os mapas de estrelas..." run ---

$linksText = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$range3 = $d.Content
$range3.Find.ClearFormatting()
$found3 = $range3.Find.Execute($linksText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $range3.Style = "GaNLinks"
}

Write-Output "Done"
